$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 535.8148
$ws.Range("I33").Value = 408.05264
$ws.Range("J33").Value = 839.25
$ws.Range("K33").Value = 408.05264
$ws.Range("L33").Value = 839.25
$ws.Range("M33").Value = -179.05264
$ws.Range("N33").Value = -1297.25
# Row 116
$ws.Range("H116").Value = 3911.0168
$ws.Range("I116").Value = 3853.0264
$ws.Range("J116").Value = 4015.9524
$ws.Range("K116").Value = 3853.0264
$ws.Range("L116").Value = 4015.9524
$ws.Range("M116").Value = -411.0264000000002
$ws.Range("N116").Value = -10899.9524
# Row 132
$ws.Range("H132").Value = 4764663.5
$ws.Range("I132").Value = 5557998.5
$ws.Range("J132").Value = 4651.8335
$ws.Range("K132").Value = 16673995.5
$ws.Range("L132").Value = 13955.5005
$ws.Range("M132").Value = -16671465.5
$ws.Range("N132").Value = -19015.5005
# Row 136
$ws.Range("H136").Value = 31035
$ws.Range("J136").Value = 31035
$ws.Range("L136").Value = 31035
$ws.Range("N136").Value = -41235
# Row 138
$ws.Range("H138").Value = 1862.52
$ws.Range("I138").Value = 618.619
$ws.Range("J138").Value = 2763.276
$ws.Range("K138").Value = 1855.857
$ws.Range("L138").Value = 8289.828
$ws.Range("M138").Value = 3284.143
$ws.Range("N138").Value = -18569.828
# Row 141
$ws.Range("H141").Value = 286304.28
$ws.Range("I141").Value = 1149.8788
$ws.Range("J141").Value = 1854653.5
$ws.Range("K141").Value = 3449.6364
$ws.Range("L141").Value = 5563960.5
$ws.Range("M141").Value = 1730.3636
$ws.Range("N141").Value = -5574320.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1238.52
$ws.Range("I32").Value = 1068.2333
$ws.Range("K32").Value = 1068.2333
$ws.Range("M32").Value = -781.2333000000001
# Row 61
$ws.Range("H61").Value = 4357.143
$ws.Range("I61").Value = 1200
$ws.Range("J61").Value = 6111.1113
$ws.Range("K61").Value = 1200
$ws.Range("L61").Value = 6111.1113
$ws.Range("M61").Value = -988
$ws.Range("N61").Value = -6535.1113
# Row 102
$ws.Range("H102").Value = 2500
$ws.Range("I102").Value = 2500
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2500
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -878
# Row 132
$ws.Range("H132").Value = 2191.9607
$ws.Range("I132").Value = 1584.8235
$ws.Range("K132").Value = 4754.470499999999
$ws.Range("M132").Value = -2224.470499999999
# Row 136
$ws.Range("H136").Value = 4357.143
$ws.Range("I136").Value = 1200
$ws.Range("J136").Value = 6111.1113
$ws.Range("K136").Value = 3600
$ws.Range("L136").Value = 18333.3339
$ws.Range("M136").Value = -1050
$ws.Range("N136").Value = -23433.3339

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 1498.7805
$ws.Range("I105").Value = 1309.56
$ws.Range("K105").Value = 1309.56
$ws.Range("M105").Value = 437.4400000000001
# Row 134
$ws.Range("H134").Value = 1845.6666
$ws.Range("I134").Value = 1305.7561
$ws.Range("J134").Value = 4059.3
$ws.Range("K134").Value = 3917.2683
$ws.Range("L134").Value = 12177.9
$ws.Range("M134").Value = -1382.2683
$ws.Range("N134").Value = -17247.9

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 1942.279
$ws.Range("I132").Value = 1407.303
$ws.Range("K132").Value = 4221.909000000001
$ws.Range("M132").Value = -1691.909000000001
# Row 133
$ws.Range("H133").Value = 21708
$ws.Range("J133").Value = 21708
$ws.Range("L133").Value = 21708
$ws.Range("N133").Value = -26768
# Row 134
$ws.Range("H134").Value = 2418.3684
$ws.Range("I134").Value = 909.6429000000001
$ws.Range("K134").Value = 2728.9287
$ws.Range("M134").Value = -193.9287000000004

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 118
$ws.Range("H118").Value = 2450
$ws.Range("J118").Value = 3000
$ws.Range("L118").Value = 9000
$ws.Range("N118").Value = -11486
# Row 131
$ws.Range("H131").Value = 1384.8462
$ws.Range("J131").Value = 1188.973
$ws.Range("L131").Value = 3566.919
$ws.Range("N131").Value = -13646.919
# Row 132
$ws.Range("H132").Value = 3150
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 4800
$ws.Range("K132").Value = 13500
$ws.Range("L132").Value = 43200
$ws.Range("M132").Value = -10970
$ws.Range("N132").Value = -48260

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 44
$ws.Range("H44").Value = 80031
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 71429940
$ws.Range("I22").Value = 142857780
$ws.Range("J22").Value = 2093.2856
$ws.Range("K22").Value = 142857780
$ws.Range("L22").Value = 2093.2856
$ws.Range("M22").Value = -142857485
$ws.Range("N22").Value = -2683.2856
# Row 27
$ws.Range("H27").Value = 71429940
$ws.Range("I27").Value = 142857780
$ws.Range("J27").Value = 2093.2856
$ws.Range("K27").Value = 142857780
$ws.Range("L27").Value = 2093.2856
$ws.Range("M27").Value = -142857673
$ws.Range("N27").Value = -2307.2856
# Row 68
$ws.Range("H68").Value = 2159.1333
$ws.Range("I68").Value = 1067.3182
$ws.Range("J68").Value = 5161.625
$ws.Range("K68").Value = 1067.3182
$ws.Range("L68").Value = 5161.625
$ws.Range("M68").Value = -318.3181999999999
$ws.Range("N68").Value = -6659.625
# Row 71
$ws.Range("H71").Value = 2159.1333
$ws.Range("I71").Value = 1067.3182
$ws.Range("J71").Value = 5161.625
$ws.Range("K71").Value = 5336.590999999999
$ws.Range("L71").Value = 25808.125
$ws.Range("M71").Value = -1592.590999999999
$ws.Range("N71").Value = -33296.125
# Row 82
$ws.Range("H82").Value = 2842.2222
$ws.Range("I82").Value = 1895
$ws.Range("J82").Value = 3600
$ws.Range("K82").Value = 1895
$ws.Range("L82").Value = 3600
$ws.Range("M82").Value = -1534
$ws.Range("N82").Value = -4322
# Row 85
$ws.Range("H85").Value = 2842.2222
$ws.Range("I85").Value = 1895
$ws.Range("J85").Value = 3600
$ws.Range("K85").Value = 1895
$ws.Range("L85").Value = 3600
$ws.Range("M85").Value = -647
$ws.Range("N85").Value = -6096
# Row 131
$ws.Range("H131").Value = 25000
$ws.Range("J131").Value = 25000
$ws.Range("L131").Value = 25000
$ws.Range("N131").Value = -35080
# Row 132
$ws.Range("H132").Value = 1821.04
$ws.Range("I132").Value = 1078.85
$ws.Range("J132").Value = 4789.8
$ws.Range("K132").Value = 3236.55
$ws.Range("L132").Value = 14369.4
$ws.Range("M132").Value = -706.5499999999997
$ws.Range("N132").Value = -19429.4
# Row 133
$ws.Range("H133").Value = 30907.691
$ws.Range("J133").Value = 30907.691
$ws.Range("L133").Value = 30907.691
$ws.Range("N133").Value = -35967.691
# Row 139
$ws.Range("H139").Value = 43782.145
$ws.Range("J139").Value = 43782.145
$ws.Range("L139").Value = 43782.145
$ws.Range("N139").Value = -54062.145

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 13703.375
$ws.Range("J54").Value = 13703.375
$ws.Range("L54").Value = 13703.375
$ws.Range("N54").Value = -14743.375
# Row 108
$ws.Range("H108").Value = 33000
$ws.Range("J108").Value = 33000
$ws.Range("L108").Value = 33000
$ws.Range("N108").Value = -40680
# Row 132
$ws.Range("H132").Value = 13568.218
$ws.Range("I132").Value = 2630.4412
$ws.Range("J132").Value = 44558.582
$ws.Range("K132").Value = 7891.323600000001
$ws.Range("L132").Value = 133675.746
$ws.Range("M132").Value = -5361.323600000001
$ws.Range("N132").Value = -138735.746
# Row 135
$ws.Range("H135").Value = 40207.5
$ws.Range("J135").Value = 40207.5
$ws.Range("L135").Value = 40207.5
$ws.Range("N135").Value = -50347.5
